# Update market/profit data on the leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
# These columns (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) hold data refreshed by the
# scheduled market-data runner; this script writes the newly-fetched values in place.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 502.64706
$ws.Range("I33").Value = 335.83334
$ws.Range("J33").Value = 593.63635
$ws.Range("K33").Value = 335.83334
$ws.Range("L33").Value = 593.63635
$ws.Range("M33").Value = -106.83334
$ws.Range("N33").Value = -1051.63635
# Row 51: A Bile Business
$ws.Range("H51").Value = 7452.2
$ws.Range("I51").Value = 7381.5
$ws.Range("J51").Value = 7499.3335
$ws.Range("K51").Value = 7381.5
$ws.Range("L51").Value = 7499.3335
$ws.Range("M51").Value = -6897.5
$ws.Range("N51").Value = -8467.333500000001
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 4582.4375
$ws.Range("I62").Value = 4582.4375
$ws.Range("K62").Value = 4582.4375
$ws.Range("M62").Value = -3958.4375
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 4582.4375
$ws.Range("I65").Value = 4582.4375
$ws.Range("K65").Value = 22912.1875
$ws.Range("M65").Value = -19792.1875
# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 3793.762
$ws.Range("I86").Value = 3178.4285
$ws.Range("K86").Value = 3178.4285
$ws.Range("M86").Value = -2055.4285
# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 3793.762
$ws.Range("I89").Value = 3178.4285
$ws.Range("K89").Value = 15892.1425
$ws.Range("M89").Value = -10276.1425
# Row 109: A Time for Peace
$ws.Range("H109").Value = 39600
$ws.Range("J109").Value = 39600
$ws.Range("L109").Value = 39600
$ws.Range("N109").Value = -42374
# Row 130: Technically Still Magic
$ws.Range("H130").Value = 42140.715
$ws.Range("J130").Value = 42140.715
$ws.Range("L130").Value = 42140.715
$ws.Range("N130").Value = -52180.715
# Row 131: Mindful Study
$ws.Range("H131").Value = 10012.523
$ws.Range("I131").Value = 8516.4375
$ws.Range("J131").Value = 14800
$ws.Range("K131").Value = 25549.3125
$ws.Range("L131").Value = 44400
$ws.Range("M131").Value = -20509.3125
$ws.Range("N131").Value = -54480
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 5494.067
$ws.Range("I137").Value = 5685.273
$ws.Range("K137").Value = 17055.819
$ws.Range("M137").Value = -14505.819

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 5879.6597
$ws.Range("I32").Value = 6903.8687
$ws.Range("K32").Value = 6903.8687
$ws.Range("M32").Value = -6616.8687
# Row 44: Very Slow Array
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976
# Row 55: Employee Retention
$ws.Range("H55").Value = 60624.6
$ws.Range("J55").Value = 112561.5
$ws.Range("L55").Value = 112561.5
$ws.Range("N55").Value = -113191.5
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4899.625
$ws.Range("I61").Value = 4499.6665
$ws.Range("K61").Value = 4499.6665
$ws.Range("M61").Value = -4287.6665
# Row 97: Ore for Me
$ws.Range("H97").Value = 3127
$ws.Range("I97").Value = 1371.1111
$ws.Range("J97").Value = 7077.75
$ws.Range("K97").Value = 1371.1111
$ws.Range("L97").Value = 7077.75
$ws.Range("M97").Value = -875.1111000000001
$ws.Range("N97").Value = -8069.75
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2737.8076
$ws.Range("I132").Value = 2511
$ws.Range("K132").Value = 7533
$ws.Range("M132").Value = -5003
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4899.625
$ws.Range("I136").Value = 4499.6665
$ws.Range("K136").Value = 13498.9995
$ws.Range("M136").Value = -10948.9995

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 25002344
$ws.Range("I86").Value = 26317204
$ws.Range("K86").Value = 26317204
$ws.Range("M86").Value = -26316081
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 25002344
$ws.Range("I89").Value = 26317204
$ws.Range("K89").Value = 131586020
$ws.Range("M89").Value = -131580404
# Row 94: High Steal
$ws.Range("H94").Value = 2581.389
$ws.Range("I94").Value = 997.3077
$ws.Range("J94").Value = 6700
$ws.Range("K94").Value = 997.3077
$ws.Range("L94").Value = 6700
$ws.Range("M94").Value = -546.3077
$ws.Range("N94").Value = -7602
# Row 138: Bladewinner
$ws.Range("H138").Value = 64283.93
$ws.Range("J138").Value = 64283.93
$ws.Range("L138").Value = 64283.93
$ws.Range("N138").Value = -74563.92999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99: O Pine
$ws.Range("H99").Value = 10202336
$ws.Range("I99").Value = 1437551.9
$ws.Range("J99").Value = 40002600
$ws.Range("K99").Value = 1437551.9
$ws.Range("L99").Value = 40002600
$ws.Range("M99").Value = -1436053.9
$ws.Range("N99").Value = -40005596
# Row 126: A Better Conductor
$ws.Range("H126").Value = 10202336
$ws.Range("I126").Value = 1437551.9
$ws.Range("J126").Value = 40002600
$ws.Range("K126").Value = 4312655.699999999
$ws.Range("L126").Value = 120007800
$ws.Range("M126").Value = -4310185.699999999
$ws.Range("N126").Value = -120012740
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1889.8096
$ws.Range("I134").Value = 1864.3889
$ws.Range("K134").Value = 5593.1667
$ws.Range("M134").Value = -3058.1667

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 1036.6
$ws.Range("I5").Value = 595
$ws.Range("J5").Value = 1856.7142
$ws.Range("K5").Value = 1785
$ws.Range("L5").Value = 5570.142599999999
$ws.Range("M5").Value = -1673
$ws.Range("N5").Value = -5794.142599999999
# Row 44: No More Dumpster Diving
$ws.Range("H44").Value = 66666730
$ws.Range("I44").Value = 83
$ws.Range("K44").Value = 249
$ws.Range("M44").Value = 149
# Row 60: Drinking to Your Health
$ws.Range("H60").Value = 823.25
$ws.Range("I60").Value = 194.81818
$ws.Range("J60").Value = 2205.8
$ws.Range("K60").Value = 584.4545400000001
$ws.Range("L60").Value = 6617.400000000001
$ws.Range("M60").Value = -333.4545400000001
$ws.Range("N60").Value = -7119.400000000001
# Row 75: Breakfast of Champions
$ws.Range("H75").Value = 2406.25
$ws.Range("I75").Value = 1237.6666
$ws.Range("J75").Value = 3107.4
$ws.Range("K75").Value = 3712.9998
$ws.Range("L75").Value = 9322.200000000001
$ws.Range("M75").Value = -2714.9998
$ws.Range("N75").Value = -11318.2
# Row 78: Emerald Soup for the Soul (L)
$ws.Range("H78").Value = 2406.25
$ws.Range("I78").Value = 1237.6666
$ws.Range("J78").Value = 3107.4
$ws.Range("K78").Value = 11138.9994
$ws.Range("L78").Value = 27966.6
$ws.Range("M78").Value = -6146.999400000001
$ws.Range("N78").Value = -37950.60000000001
# Row 94: All You Can Stomach
$ws.Range("H94").Value = 3519
$ws.Range("I94").Value = 2395
$ws.Range("J94").Value = 3800
$ws.Range("K94").Value = 7185
$ws.Range("L94").Value = 11400
$ws.Range("M94").Value = -6509
$ws.Range("N94").Value = -12752
# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 1491.3334
$ws.Range("J121").Value = 1749.5
$ws.Range("L121").Value = 5248.5
$ws.Range("N121").Value = -7868.5
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1036.6
$ws.Range("I135").Value = 595
$ws.Range("J135").Value = 1856.7142
$ws.Range("K135").Value = 5355
$ws.Range("L135").Value = 16710.4278
$ws.Range("M135").Value = -2820
$ws.Range("N135").Value = -21780.4278

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 44: Actually, It's Loyalty
$ws.Range("H44").Value = 19666.5
$ws.Range("I44").Value = 19500
$ws.Range("J44").Value = 19999.5
$ws.Range("K44").Value = 19500
$ws.Range("L44").Value = 19999.5
$ws.Range("M44").Value = -18904
$ws.Range("N44").Value = -21191.5
# Row 132: On Board for Lar
$ws.Range("H132").Value = 7259.6284
$ws.Range("I132").Value = 6235.207
$ws.Range("J132").Value = 12211
$ws.Range("K132").Value = 18705.621
$ws.Range("L132").Value = 36633
$ws.Range("M132").Value = -16175.621
$ws.Range("N132").Value = -41693

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 69: Maybe He's a Lion
$ws.Range("H69").Value = 67581.5
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14189
# Row 72: The Wyvern of It (L)
$ws.Range("H72").Value = 67581.5
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40944
# Row 122: Hell on Leather
$ws.Range("H122").Value = 3555.875
$ws.Range("J122").Value = 3989
$ws.Range("L122").Value = 11967
$ws.Range("N122").Value = -16867

